# ------------------------------------------------------------------
# edit.ps1 - applies "feat: add 2022-Q3 data" to 300308-中际旭创.xlsx
#
# Summary of the change:
#  1. "总计" (sheet1) gets a new row 2 for "2022-Q3" (count=17, value=3.73);
#     all other rows shift down by one and the running index in column A
#     is renumbered 0..7.
#  2. A brand-new worksheet named "2022-Q3" is inserted right after "总计"
#     (i.e. it becomes the second tab), holding the per-fund holdings
#     table (same layout as the neighbouring "2022-Q2" sheet).
#  3. Every other quarter tab keeps its data untouched and simply shifts
#     one position to the right in the tab strip.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert the 2022-Q3 row at the top of the
#    data block (row 2), push the rest down.
# ------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Shift existing data rows 2..8 down to 3..9 (bottom-up so sources are
# never clobbered before they are read). Range.Copy brings the value
# AND the style along, so column A keeps its "s=2" formatting intact.
for ($r = 8; $r -ge 2; $r--) {
    $nr = $r + 1
    $totalSheet.Range("A" + $r + ":D" + $r).Copy($totalSheet.Range("A" + $nr))
}

# Write the new 2022-Q3 figures into row 2.
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 17
$totalSheet.Range("D2").Value = 3.73

# Renumber the running index in column A (0-based) for every data row.
for ($r = 2; $r -le 9; $r++) {
    $totalSheet.Range("A" + $r).Value = $r - 2
}

# ------------------------------------------------------------------
# 2) Insert the new "2022-Q3" worksheet as the 2nd tab (right after
#    "总计", right before the existing "2022-Q2" tab).
# ------------------------------------------------------------------
$beforeTarget = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($beforeTarget)
$newSheet.Name = "2022-Q3"

# NOTE: Worksheets.Add() shifts every sheet's position, which can
# invalidate sheet references captured beforehand — so re-fetch the
# "2022-Q2" sheet fresh, *after* the insert, before using it as a
# formatting/layout template.
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# Bring over the header row + per-row formatting by copying the shape
# of the neighbouring "2022-Q2" sheet (same columns/styles), then
# overwrite every cell with the real 2022-Q3 values below.
$q2Sheet.Range("B1:H1").Copy($newSheet.Range("B1"))
$q2Sheet.Range("A2:H18").Copy($newSheet.Range("A2"))

# Columns B-G hold fund codes / figures that must stay TEXT (several
# fund codes have leading zeros, e.g. "009686"), so force a text
# number format before writing the values — otherwise Excel would
# silently coerce them to numbers and drop the leading zeros.
$newSheet.Range("B2:G18").NumberFormat = "@"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "501022"
$newSheet.Range("C2").Value = "银华鑫盛灵活配置混合（LOF）A"
$newSheet.Range("D2").Value = "41.87"
$newSheet.Range("E2").Value = "72.03"
$newSheet.Range("F2").Value = "2.49"
$newSheet.Range("G2").Value = "1.0426"
$newSheet.Range("H2").Value = 7
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "161834"
$newSheet.Range("C3").Value = "银华鑫锐灵活配置混合（LOF）A"
$newSheet.Range("D3").Value = "36.66"
$newSheet.Range("E3").Value = "75.53"
$newSheet.Range("F3").Value = "2.62"
$newSheet.Range("G3").Value = "0.9605"
$newSheet.Range("H3").Value = 7
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "009686"
$newSheet.Range("C4").Value = "华夏磐利一年定期开放混合A"
$newSheet.Range("D4").Value = "10.76"
$newSheet.Range("E4").Value = "64.78"
$newSheet.Range("F4").Value = "2.71"
$newSheet.Range("G4").Value = "0.2916"
$newSheet.Range("H4").Value = 5
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "240008"
$newSheet.Range("C5").Value = "华宝收益增长混合A"
$newSheet.Range("D5").Value = "8.39"
$newSheet.Range("E5").Value = "93.74"
$newSheet.Range("F5").Value = "3.44"
$newSheet.Range("G5").Value = "0.2886"
$newSheet.Range("H5").Value = 9
$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "012370"
$newSheet.Range("C6").Value = "银华鑫利一年持有期混合"
$newSheet.Range("D6").Value = "9.61"
$newSheet.Range("E6").Value = "73.11"
$newSheet.Range("F6").Value = "2.15"
$newSheet.Range("G6").Value = "0.2066"
$newSheet.Range("H6").Value = 8
$newSheet.Range("A7").Value = 5
$newSheet.Range("B7").Value = "014048"
$newSheet.Range("C7").Value = "银华鑫盛灵活配置混合（LOF）C"
$newSheet.Range("D7").Value = "8.19"
$newSheet.Range("E7").Value = "72.03"
$newSheet.Range("F7").Value = "2.49"
$newSheet.Range("G7").Value = "0.2039"
$newSheet.Range("H7").Value = 7
$newSheet.Range("A8").Value = 6
$newSheet.Range("B8").Value = "002770"
$newSheet.Range("C8").Value = "安信新回报灵活配置混合A"
$newSheet.Range("D8").Value = "4.43"
$newSheet.Range("E8").Value = "88.08"
$newSheet.Range("F8").Value = "3.88"
$newSheet.Range("G8").Value = "0.1719"
$newSheet.Range("H8").Value = 8
$newSheet.Range("A9").Value = 7
$newSheet.Range("B9").Value = "002771"
$newSheet.Range("C9").Value = "安信新回报灵活配置混合C"
$newSheet.Range("D9").Value = "3.64"
$newSheet.Range("E9").Value = "88.08"
$newSheet.Range("F9").Value = "3.88"
$newSheet.Range("G9").Value = "0.1412"
$newSheet.Range("H9").Value = 8
$newSheet.Range("A10").Value = 8
$newSheet.Range("B10").Value = "004350"
$newSheet.Range("C10").Value = "汇丰晋信价值先锋股票A"
$newSheet.Range("D10").Value = "5.31"
$newSheet.Range("E10").Value = "94.44"
$newSheet.Range("F10").Value = "2.52"
$newSheet.Range("G10").Value = "0.1338"
$newSheet.Range("H10").Value = 10
$newSheet.Range("A11").Value = 9
$newSheet.Range("B11").Value = "014349"
$newSheet.Range("C11").Value = "银华鑫锐灵活配置混合（LOF）C"
$newSheet.Range("D11").Value = "4.69"
$newSheet.Range("E11").Value = "75.53"
$newSheet.Range("F11").Value = "2.62"
$newSheet.Range("G11").Value = "0.1229"
$newSheet.Range("H11").Value = 7
$newSheet.Range("A12").Value = 10
$newSheet.Range("B12").Value = "010033"
$newSheet.Range("C12").Value = "安信成长精选混合A"
$newSheet.Range("D12").Value = "1.66"
$newSheet.Range("E12").Value = "88.80"
$newSheet.Range("F12").Value = "4.12"
$newSheet.Range("G12").Value = "0.0684"
$newSheet.Range("H12").Value = 8
$newSheet.Range("A13").Value = 11
$newSheet.Range("B13").Value = "012321"
$newSheet.Range("C13").Value = "东财中证云计算指数增强A"
$newSheet.Range("D13").Value = "1.31"
$newSheet.Range("E13").Value = "93.47"
$newSheet.Range("F13").Value = "4.07"
$newSheet.Range("G13").Value = "0.0533"
$newSheet.Range("H13").Value = 8
$newSheet.Range("A14").Value = 12
$newSheet.Range("B14").Value = "012322"
$newSheet.Range("C14").Value = "东财中证云计算指数增强C"
$newSheet.Range("D14").Value = "0.51"
$newSheet.Range("E14").Value = "93.47"
$newSheet.Range("F14").Value = "4.07"
$newSheet.Range("G14").Value = "0.0208"
$newSheet.Range("H14").Value = 8
$newSheet.Range("A15").Value = 13
$newSheet.Range("B15").Value = "009687"
$newSheet.Range("C15").Value = "华夏磐利一年定期开放混合C"
$newSheet.Range("D15").Value = "0.43"
$newSheet.Range("E15").Value = "64.78"
$newSheet.Range("F15").Value = "2.71"
$newSheet.Range("G15").Value = "0.0117"
$newSheet.Range("H15").Value = 5
$newSheet.Range("A16").Value = 14
$newSheet.Range("B16").Value = "010034"
$newSheet.Range("C16").Value = "安信成长精选混合C"
$newSheet.Range("D16").Value = "0.24"
$newSheet.Range("E16").Value = "88.80"
$newSheet.Range("F16").Value = "4.12"
$newSheet.Range("G16").Value = "0.0099"
$newSheet.Range("H16").Value = 8
$newSheet.Range("A17").Value = 15
$newSheet.Range("B17").Value = "015573"
$newSheet.Range("C17").Value = "华宝收益增长混合C"
$newSheet.Range("D17").Value = "0.04"
$newSheet.Range("E17").Value = "93.74"
$newSheet.Range("F17").Value = "3.44"
$newSheet.Range("G17").Value = "0.0014"
$newSheet.Range("H17").Value = 9
$newSheet.Range("A18").Value = 16
$newSheet.Range("B18").Value = "015364"
$newSheet.Range("C18").Value = "汇丰晋信价值先锋股票C"
$newSheet.Range("D18").Value = "0.02"
$newSheet.Range("E18").Value = "94.44"
$newSheet.Range("F18").Value = "2.52"
$newSheet.Range("G18").Value = "0.0005"
$newSheet.Range("H18").Value = 10

